$d = $word.ActiveDocument

$d.Content.Find.Execute("Rationalizing the denominator (priority)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Rationalizing the denominator (priority, Max)", 2)

$d.Content.Find.Execute("Artihmetic on complex numbers", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Artihmetic on complex numbers (Charlotte)", 2)

$d.Content.Find.Execute("Trigonometry and integration", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Trigonometry and integration (Sophie C)", 2)

$d.Content.Find.Execute("Solving simultaneous equations (priority)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Solving simultaneous equations (priority, Ollie)", 2)

$d.Content.Find.Execute("Hypothesis testing", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Introduction to hypothesis testing (Ellie)", 2)
